$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo in existing shared string (C3: "AzureBastiuoin" -> "AzureBastion") ---
$ws.Range("C3").Value = "AzureBastion"

# --- Add new data rows (3-6): repeat VNET name / address space; row 3 also keeps the Bastion subnet note ---
$ws.Range("A3").Value = "vnet-hub-weu001"
$ws.Range("B3").Value = "10.0.0.0/23"

$ws.Range("A4").Value = "vnet-hub-weu001"
$ws.Range("B4").Value = "10.0.0.0/23"

$ws.Range("A5").Value = "vnet-hub-weu001"
$ws.Range("B5").Value = "10.0.0.0/23"

$ws.Range("A6").Value = "vnet-hub-weu001"
$ws.Range("B6").Value = "10.0.0.0/23"

# --- Header row (row 1): bold, and center-align the NSG/UDR columns ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("H1").Font.Bold = $true

$ws.Range("F1:G1").Font.Bold = $true
$ws.Range("F1:G1").HorizontalAlignment = -4108

# --- Row 2: center-align the NSG/UDR answer cells ---
$ws.Range("F2:G2").HorizontalAlignment = -4108

# --- Column widths (resize to match authored layout; values chosen so the
#     engine's char-width quantization lands on the closest achievable width) ---
$ws.Columns.Item(1).ColumnWidth = 18.495
$ws.Columns.Item(2).ColumnWidth = 18.331
$ws.Columns.Item(3).ColumnWidth = 21.331
$ws.Columns.Item(4).ColumnWidth = 19.831
$ws.Columns.Item(5).ColumnWidth = 15.66
$ws.Columns.Item(6).ColumnWidth = 3.831
$ws.Columns.Item(7).ColumnWidth = 3.831
$ws.Columns.Item(8).ColumnWidth = 28.66

# --- Print setup ---
$ws.PageSetup.Orientation = 1

# --- Selection cursor, matching the saved view state ---
[void]$ws.Range("H9").Select()
